# Refresh the cryptos list with the latest GitHub Actions scrape.
# Each cell is written with its exact final text -- values such as
# "1.00" or "0.611" are coin prices stored as text in this sheet, so
# a leading quote is used to keep Excel from re-interpreting them as
# numbers (matching the existing text-formatted Price column).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '67.426.65'
$ws.Range("E2").Value = '  -0.89%  '

# Row 3
$ws.Range("D3").Value = '3.229.46'
$ws.Range("E3").Value = '  -1.25%  '

# Row 4
$ws.Range("E4").Value = '  +0.04%  '

# Row 5
$ws.Range("D5").Value = "'578.36"
$ws.Range("E5").Value = '  -1.64%  '

# Row 6
$ws.Range("D6").Value = "'184.25"
$ws.Range("E6").Value = '  -1.40%  '

# Row 7
$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").Value = "'0.611"
$ws.Range("E7").Value = '  +1.76%  '

# Row 8
$ws.Range("B8").Value = 'USDC'
$ws.Range("C8").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = '  +0.05%  '

# Row 9
$ws.Range("D9").Value = '3.227.58'

# Row 10
$ws.Range("E10").Value = '  -3.33%  '

# Row 11
$ws.Range("D11").Value = "'6.58"
$ws.Range("E11").Value = '  -2.20%  '

# Row 12
$ws.Range("D12").Value = "'0.411"
$ws.Range("E12").Value = '  -1.43%  '

# Row 13
$ws.Range("D13").Value = '3.784.50'
$ws.Range("E13").Value = '  -1.33%  '

# Row 15
$ws.Range("D15").Value = "'27.64"
$ws.Range("E15").Value = '  -3.25%  '

# Row 16
$ws.Range("D16").Value = '67.479.64'
$ws.Range("E16").Value = '  -0.77%  '

# Row 17
$ws.Range("D17").Value = "'0.0000169"
$ws.Range("E17").Value = '  -2.05%  '

# Row 18
$ws.Range("D18").Value = '3.231.67'
$ws.Range("E18").Value = '  -1.08%  '

# Row 19
$ws.Range("D19").Value = "'5.75"

# Row 20
$ws.Range("D20").Value = "'13.45"
$ws.Range("E20").Value = '  -1.19%  '

# Row 21
$ws.Range("D21").Value = "'395.61"
$ws.Range("E21").Value = '  +3.89%  '

# Row 22
$ws.Range("D22").Value = "'7.55"
$ws.Range("E22").Value = '  -2.31%  '

# Row 23
$ws.Range("E23").Value = '  +0.21%  '

# Row 24
$ws.Range("D24").Value = "'71.20"
$ws.Range("E24").Value = '  -0.37%  '

# Row 25
$ws.Range("D25").Value = "'0.516"
$ws.Range("E25").Value = '  +0.25%  '

# Row 26
$ws.Range("D26").Value = "'0.0000118"
$ws.Range("E26").Value = '  -2.45%  '

# Row 27
$ws.Range("D27").Value = "'0.187"
$ws.Range("E27").Value = '  -1.02%  '

# Row 28
$ws.Range("D28").Value = "'9.53"
$ws.Range("E28").Value = '  -2.86%  '

# Row 29
$ws.Range("E29").Value = '  +0.27%  '

# Row 30
$ws.Range("E30").Value = '  -2.19%  '

# Row 31
$ws.Range("D31").Value = "'5.55"
$ws.Range("E31").Value = '  -4.18%  '

# Row 32
$ws.Range("D32").Value = "'22.56"
$ws.Range("E32").Value = '  -1.31%  '

# Row 33
$ws.Range("E33").Value = '  -3.15%  '

# Row 34
$ws.Range("E34").Value = '  -1.87%  '

# Row 35
$ws.Range("E35").Value = '  +0.03%  '

# Row 36
$ws.Range("D36").Value = "'160.76"
$ws.Range("E36").Value = '  -1.15%  '

# Row 37
$ws.Range("E37").Value = '  -4.37%  '

# Row 38
$ws.Range("E38").Value = '  +0.94%  '

# Row 39
$ws.Range("D39").Value = "'26.42"
$ws.Range("E39").Value = '  -0.75%  '

# Row 40
$ws.Range("E40").Value = '  -4.33%  '

# Row 41
$ws.Range("D41").Value = "'4.54"
$ws.Range("E41").Value = '  -1.17%  '

# Row 42
$ws.Range("D42").Value = "'6.50"
$ws.Range("E42").Value = '  -4.52%  '

# Row 43
$ws.Range("D43").Value = "'2.47"
$ws.Range("E43").Value = '  -5.90%  '

# Row 44
$ws.Range("E44").Value = '  -0.86%  '

# Row 45
$ws.Range("D45").Value = "'40.47"
$ws.Range("E45").Value = '  -1.77%  '

# Row 46
$ws.Range("D46").Value = '2.600.43'
$ws.Range("E46").Value = '  -1.68%  '

# Row 47
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").Value = "'24.57"
$ws.Range("E47").Value = '  -3.38%  '

# Row 48
$ws.Range("B48").Value = 'Bittensor'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D48").Value = "'334.33"
$ws.Range("E48").Value = '  -2.55%  '

# Row 49
$ws.Range("D49").Value = "'0.0278"
$ws.Range("E49").Value = '  -2.39%  '

# Row 50
$ws.Range("D50").Value = "'6.27"
$ws.Range("E50").Value = '  +0.15%  '

# Row 51
$ws.Range("D51").Value = "'0.102"
$ws.Range("E51").Value = '  -0.45%  '
